$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.906.43"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.604.48"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'210.60"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.484"
$ws.Range("E7").Value = "  -3.92%  "
$ws.Range("D8").Value = "'0.0615"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "'17.97"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "1.828.20"
$ws.Range("D13").Value = "1.618.08"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  -3.24%  "
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").Value = "25.898.34"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'61.22"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'189.90"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'9.37"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("D25").Value = "'142.27"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("D28").Value = "'6.56"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").Value = "'15.00"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").Value = "'3.04"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").Value = "1.117.04"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "'0.811"
$ws.Range("E38").Value = "  -6.79%  "
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").Value = "'96.01"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "1.740.62"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "'0.745"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("D44").Value = "'5.06"
$ws.Range("E44").Value = "  -4.06%  "
$ws.Range("D45").Value = "0.0₆0113"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'53.41"
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D51").Value = "'7.40"
$ws.Range("E51").Value = "  -1.73%  "
